$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.895.17'
$ws.Range("E2").Value = '  -0.55%  '

# Row 3
$ws.Range("D3").Value = '1.641.81'
$ws.Range("E3").Value = '  -0.04%  '

# Row 4
$ws.Range("D4").Value = '''1.010'
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").Value = '''215.54'
$ws.Range("E5").Value = '  -0.49%  '

# Row 6
$ws.Range("D6").Value = '''0.5063'
$ws.Range("E6").Value = '  +0.16%  '

# Row 7
$ws.Range("D7").Value = '''1.006'
$ws.Range("E7").Value = '  -0.29%  '

# Row 8
$ws.Range("D8").Value = '''0.2582'
$ws.Range("E8").Value = '  +0.22%  '

# Row 9
$ws.Range("D9").Value = '''0.06419'
$ws.Range("E9").Value = '  -0.51%  '

# Row 10
$ws.Range("D10").Value = '''19.79'
$ws.Range("E10").Value = '  +0.83%  '

# Row 11
$ws.Range("D11").Value = '''0.07810'
$ws.Range("E11").Value = '  +1.39%  '

# Row 12
$ws.Range("D12").Value = '1.663.80'
$ws.Range("E12").Value = '  +0.96%  '

# Row 13
$ws.Range("D13").Value = '''4.307'
$ws.Range("E13").Value = '  +1.25%  '

# Row 14
$ws.Range("D14").Value = '''0.5462'
$ws.Range("E14").Value = '  -0.20%  '

# Row 15
$ws.Range("D15").Value = '0.0₅7899'
$ws.Range("E15").Value = '  -0.76%  '

# Row 16
$ws.Range("D16").Value = '''65.07'
$ws.Range("E16").Value = '  +2.22%  '

# Row 17
$ws.Range("D17").Value = '25.964.33'
$ws.Range("E17").Value = '  -0.23%  '

# Row 18
$ws.Range("D18").Value = '''1.008'
$ws.Range("E18").Value = '  -0.15%  '

# Row 19
$ws.Range("D19").Value = '''199.06'
$ws.Range("E19").Value = '  -2.44%  '

# Row 20
$ws.Range("D20").Value = '''4.419'
$ws.Range("E20").Value = '  +2.35%  '

# Row 21
$ws.Range("D21").Value = '''10.02'
$ws.Range("E21").Value = '  -0.09%  '

# Row 22
$ws.Range("D22").Value = '''6.025'
$ws.Range("E22").Value = '  +0.31%  '

# Row 23
$ws.Range("D23").Value = '''1.008'
$ws.Range("E23").Value = '  -0.20%  '

# Row 24
$ws.Range("D24").Value = '''1.900'
$ws.Range("E24").Value = '  -2.94%  '

# Row 25
$ws.Range("D25").Value = '''140.44'
$ws.Range("E25").Value = '  -1.23%  '

# Row 26
$ws.Range("D26").Value = '''0.1146'
$ws.Range("E26").Value = '  -0.80%  '

# Row 27
$ws.Range("D27").Value = '''6.908'
$ws.Range("E27").Value = '  +2.46%  '

# Row 28
$ws.Range("D28").Value = '''15.76'
$ws.Range("E28").Value = '  +0.06%  '

# Row 29
$ws.Range("D29").Value = '''1.242'
$ws.Range("E29").Value = '  -0.29%  '

# Row 30
$ws.Range("D30").Value = '''0.05048'
$ws.Range("E30").Value = '  -0.76%  '

# Row 31
$ws.Range("D31").Value = '''3.270'
$ws.Range("E31").Value = '  -0.04%  '

# Row 32
$ws.Range("D32").Value = '''3.217'
$ws.Range("E32").Value = '  +0.61%  '

# Row 33
$ws.Range("D33").Value = '''1.535'
$ws.Range("E33").Value = '  -0.62%  '

# Row 34
$ws.Range("D34").Value = '''2.371'
$ws.Range("E34").Value = '  +0.41%  '

# Row 35
$ws.Range("D35").Value = '''0.8974'
$ws.Range("E35").Value = '  -0.04%  '

# Row 36
$ws.Range("D36").Value = '''2.598'
$ws.Range("E36").Value = '  -1.38%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.5567'
$ws.Range("E37").Value = '  -0.55%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.137.61'
$ws.Range("E38").Value = '  -2.97%  '

# Row 39
$ws.Range("D39").Value = '''0.01559'
$ws.Range("E39").Value = '  -0.65%  '

# Row 40
$ws.Range("E40").Value = '  -0.21%  '

# Row 41
$ws.Range("D41").Value = '''5.732'
$ws.Range("E41").Value = '  +0.93%  '

# Row 42
$ws.Range("D42").Value = '''0.8191'
$ws.Range("E42").Value = '  +1.19%  '

# Row 43
$ws.Range("D43").Value = '''99.98'
$ws.Range("E43").Value = '  -0.09%  '

# Row 44
$ws.Range("D44").Value = '0.0₈121'
$ws.Range("E44").Value = '  +8.48%  '

# Row 45
$ws.Range("D45").Value = '1.777.49'
$ws.Range("E45").Value = '  -0.22%  '

# Row 46
$ws.Range("D46").Value = '''0.4550'
$ws.Range("E46").Value = '  +0.29%  '

# Row 47
$ws.Range("D47").Value = '''55.54'
$ws.Range("E47").Value = '  +0.63%  '

# Row 48
$ws.Range("D48").Value = '''1.007'
$ws.Range("E48").Value = '  -0.19%  '

# Row 49
$ws.Range("D49").Value = '''0.05095'
$ws.Range("E49").Value = '  +1.03%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.09559'
$ws.Range("E50").Value = '  +2.76%  '

# Row 51
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '''1.006'
$ws.Range("E51").Value = '  -0.25%  '
